$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Before: last data row is row 38 (highlighted "last row" formatting,
#         B38 style s=6, C38 style s=7), merged cell C37:C38.
# After:  a new row 39 is appended (date 45737, hours 6) and becomes
#         the new "last row"; row 38 becomes a normal data row; the
#         merged note cell grows to C37:C39.
# ------------------------------------------------------------------

# 1) Carry the current "last row" look (row 38) down onto the new row 39
#    before anything else changes.
$ws.Range("A38:C38").Copy()
$ws.Range("A39:C39").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# 2) Row 38 becomes an ordinary row now, so restyle it like row 37.
$ws.Range("A37:B37").Copy()
$ws.Range("A38:B38").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# 3) Grow the merged note cell from C37:C38 to C37:C39.
$ws.Range("C37:C39").Merge()

# 4) The merge operation above re-derives formatting for the member
#    cells it touches; reapply the note column's normal style on top
#    so C38/C39 match the rest of the note column.
$ws.Range("C37").Copy()
$ws.Range("C38").PasteSpecial(-4122)       # xlPasteFormats
$ws.Range("C39").PasteSpecial(-4122)       # xlPasteFormats
$excel.CutCopyMode = $false

# 5) Fill in the new row's data.
$ws.Range("A39").Value = 45737
$ws.Range("B39").Value = 6

# 6) Match the selection state captured in the saved workbook.
$ws.Range("C37:C39").Select()
